$d = $word.ActiveDocument

# Find the "Features " paragraph and insert a new paragraph right after it
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -match "^Features\s*\r?$") {
        $p.Range.InsertParagraphAfter()
        break
    }
}

# Now find the newly-created (empty) paragraph that follows "Features " and set its text
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -match "^Features\s*\r?$") {
        $newPara = $p.Next()
        $newPara.Range.Text = "Following are the features"
        break
    }
}
